$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "walkingToRunning"
$ws.Range("C2").Value = 0.2934486865997314
$ws.Range("D2").Value = -24.41073989868164
$ws.Range("E2").Value = -1.238773345947266
$ws.Range("F2").Value = -1.428148408217626
$ws.Range("G2").Value = 1.415199639612464
$ws.Range("H2").Value = 1.873757413392585

$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "walkingToRunning"
$ws.Range("C3").Value = 2.638106346130371
$ws.Range("D3").Value = -17.53305244445801
$ws.Range("E3").Value = -1.380356788635254
$ws.Range("F3").Value = -0.3362366334978593
$ws.Range("G3").Value = 0.3509886704687749
$ws.Range("H3").Value = 0.4294003237217614

$ws.Range("A4").Value = 200
$ws.Range("B4").Value = "walkingToRunning"
$ws.Range("C4").Value = 2.457366943359375
$ws.Range("D4").Value = -12.19800853729248
$ws.Range("E4").Value = -2.702709674835205
$ws.Range("F4").Value = 1.29839596625184
$ws.Range("G4").Value = 0.07410049284516276
$ws.Range("H4").Value = -0.7659658990882879

$ws.Range("A5").Value = 300
$ws.Range("B5").Value = "walkingToRunning"
$ws.Range("C5").Value = 1.40507447719574
$ws.Range("D5").Value = -11.89906311035156
$ws.Range("E5").Value = -5.299991607666016
$ws.Range("F5").Value = 1.272031672326387
$ws.Range("G5").Value = -1.028642801671032
$ws.Range("H5").Value = -0.5943050091996711

$ws.Range("A6").Value = 400
$ws.Range("B6").Value = "walkingToRunning"
$ws.Range("C6").Value = 7.022602081298828
$ws.Range("D6").Value = -22.80319976806641
$ws.Range("E6").Value = 12.74471092224121
$ws.Range("F6").Value = 0.3833239707559676
$ws.Range("G6").Value = -2.850368825271993
$ws.Range("H6").Value = -0.3757410840134667

$ws.Range("A7").Value = 500
$ws.Range("B7").Value = "walkingToRunning"
$ws.Range("C7").Value = -13.86397647857666
$ws.Range("D7").Value = -20.78690338134766
$ws.Range("E7").Value = 3.63317346572876
$ws.Range("F7").Value = 0.4730243273766702
$ws.Range("G7").Value = -4.100493916726186
$ws.Range("H7").Value = -1.476941543751531

$ws.Range("A8").Value = 600
$ws.Range("B8").Value = "walkingToRunning"
$ws.Range("C8").Value = -2.486610889434814
$ws.Range("D8").Value = -19.02692794799805
$ws.Range("E8").Value = -5.421220779418945
$ws.Range("F8").Value = 0.5741065289701482
$ws.Range("G8").Value = -2.503226176839088
$ws.Range("H8").Value = -1.719519167368706

$ws.Range("A9").Value = 700
$ws.Range("B9").Value = "walkingToRunning"
$ws.Range("C9").Value = 2.134828090667725
$ws.Range("D9").Value = -14.44269371032715
$ws.Range("E9").Value = -4.357526302337647
$ws.Range("F9").Value = 0.9217050730962066
$ws.Range("G9").Value = 0.2292654945401817
$ws.Range("H9").Value = -1.401101356942718

$ws.Range("A10").Value = 800
$ws.Range("B10").Value = "walkingToRunning"
$ws.Range("C10").Value = 1.494446396827698
$ws.Range("D10").Value = -15.40649795532227
$ws.Range("E10").Value = -6.447004318237305
$ws.Range("F10").Value = 1.19596545062822
$ws.Range("G10").Value = -0.2426663118644866
$ws.Range("H10").Value = -1.187804716979447

$ws.Range("A11").Value = 900
$ws.Range("B11").Value = "walkingToRunning"
$ws.Range("C11").Value = -0.8428599834442139
$ws.Range("D11").Value = -13.50097179412842
$ws.Range("E11").Value = -3.483202457427979
$ws.Range("F11").Value = 0.9760257120062067
$ws.Range("G11").Value = -0.1424154897952448
$ws.Range("H11").Value = -0.8590546718822623

$ws.Range("A12").Value = 1000
$ws.Range("B12").Value = "walkingToRunning"
$ws.Range("C12").Value = -0.948378562927246
$ws.Range("D12").Value = -17.74329948425293
$ws.Range("E12").Value = -5.102754592895508
$ws.Range("F12").Value = 0.3539368301300073
$ws.Range("G12").Value = 0.2236679532871066
$ws.Range("H12").Value = -0.08171037576295181

$ws.Range("A13").Value = 1100
$ws.Range("B13").Value = "walkingToRunning"
$ws.Range("C13").Value = 2.71714448928833
$ws.Range("D13").Value = -15.43135833740234
$ws.Range("E13").Value = -0.4548602104187011
$ws.Range("F13").Value = -1.62912566723419
$ws.Range("G13").Value = 1.116381197617942
$ws.Range("H13").Value = 1.465065406477321

$ws.Range("A14").Value = 1200
$ws.Range("B14").Value = "walkingToRunning"
$ws.Range("C14").Value = 5.847911357879639
$ws.Range("D14").Value = -17.33660507202148
$ws.Range("E14").Value = 6.749421119689941
$ws.Range("F14").Value = -4.661672166352781
$ws.Range("G14").Value = 0.5886271724621748
$ws.Range("H14").Value = 3.194017144586767

$ws.Range("A15").Value = 1300
$ws.Range("B15").Value = "walkingToRunning"
$ws.Range("C15").Value = 16.91429901123047
$ws.Range("D15").Value = -23.71953201293945
$ws.Range("E15").Value = 22.54716873168945
$ws.Range("F15").Value = -4.62152429700339
$ws.Range("G15").Value = 2.452606123751813
$ws.Range("H15").Value = 2.296193793690959

$ws.Range("A16").Value = 1400
$ws.Range("B16").Value = "walkingToRunning"
$ws.Range("C16").Value = 7.127767562866211
$ws.Range("D16").Value = -25.82431030273437
$ws.Range("E16").Value = 2.553998231887817
$ws.Range("F16").Value = 2.544737715562746
$ws.Range("G16").Value = 2.399413680678499
$ws.Range("H16").Value = -0.8650669424736119

$ws.Range("A17").Value = 1500
$ws.Range("B17").Value = "walkingToRunning"
$ws.Range("C17").Value = -12.36970138549805
$ws.Range("D17").Value = -11.04128646850586
$ws.Range("E17").Value = -19.66802978515625
$ws.Range("F17").Value = 4.244821196552569
$ws.Range("G17").Value = -2.53188033121541
$ws.Range("H17").Value = -0.9414239910695861

$ws.Range("A18").Value = 1600
$ws.Range("B18").Value = "walkingToRunning"
$ws.Range("C18").Value = -13.15851974487305
$ws.Range("D18").Value = -18.03956413269043
$ws.Range("E18").Value = -7.630072116851807
$ws.Range("F18").Value = 4.651211643570909
$ws.Range("G18").Value = 2.098859643144377
$ws.Range("H18").Value = -0.957986710256305

$ws.Range("A19").Value = 1700
$ws.Range("B19").Value = "walkingToRunning"
$ws.Range("C19").Value = 15.48420429229736
$ws.Range("D19").Value = -32.6015739440918
$ws.Range("E19").Value = 28.89531707763672
$ws.Range("F19").Value = 1.55922644780575
$ws.Range("G19").Value = -0.8341632472193089
$ws.Range("H19").Value = -1.370352138012565

$ws.Range("A20").Value = 1800
$ws.Range("B20").Value = "walkingToRunning"
$ws.Range("C20").Value = 12.63175868988037
$ws.Range("D20").Value = -10.51330184936523
$ws.Range("E20").Value = 13.76239585876465
$ws.Range("F20").Value = -7.001311344414123
$ws.Range("G20").Value = -5.33526933457129
$ws.Range("H20").Value = 2.328884744116308

$ws.Range("A21").Value = 1900
$ws.Range("B21").Value = "walkingToRunning"
$ws.Range("C21").Value = 19.2936897277832
$ws.Range("D21").Value = -10.78378295898438
$ws.Range("E21").Value = 9.518977165222168
$ws.Range("F21").Value = -9.313379083612446
$ws.Range("G21").Value = -15.71127955588033
$ws.Range("H21").Value = 2.31600133739272

$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "walkingToRunning"
$ws.Range("C22").Value = 6.209894180297852
$ws.Range("D22").Value = -38.24882507324219
$ws.Range("E22").Value = 16.8684253692627
$ws.Range("F22").Value = 5.349856077525013
$ws.Range("G22").Value = -9.149664600836861
$ws.Range("H22").Value = 4.022751534557007

$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "walkingToRunning"
$ws.Range("C23").Value = 14.28531169891357
$ws.Range("D23").Value = -6.773443222045898
$ws.Range("E23").Value = -8.362334251403809
$ws.Range("F23").Value = 4.662202395196333
$ws.Range("G23").Value = -3.257402020626782
$ws.Range("H23").Value = 0.9362263946735492

$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "walkingToRunning"
$ws.Range("C24").Value = -12.19628238677978
$ws.Range("D24").Value = -18.89439392089844
$ws.Range("E24").Value = -11.69162940979004
$ws.Range("F24").Value = 3.284184457631085
$ws.Range("G24").Value = 3.942425561567058
$ws.Range("H24").Value = 0.4433726047567848

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "walkingToRunning"
$ws.Range("C25").Value = 27.13467597961425
$ws.Range("D25").Value = -28.36122703552246
$ws.Range("E25").Value = 6.44420337677002
$ws.Range("F25").Value = -3.787294496909284
$ws.Range("G25").Value = 1.309673618566504
$ws.Range("H25").Value = 1.816239040716114

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "walkingToRunning"
$ws.Range("C26").Value = -3.487966060638428
$ws.Range("D26").Value = -1.607144355773926
$ws.Range("E26").Value = -0.4285287857055664
$ws.Range("F26").Value = -4.766657801131926
$ws.Range("G26").Value = 2.415746399837134
$ws.Range("H26").Value = 2.136085600870617

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "walkingToRunning"
$ws.Range("C27").Value = -40.79524230957031
$ws.Range("D27").Value = -87.16375732421875
$ws.Range("E27").Value = 55.55958557128906
$ws.Range("F27").Value = -3.364129057669638
$ws.Range("G27").Value = 9.653811729262191
$ws.Range("H27").Value = -2.774162345706682

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "walkingToRunning"
$ws.Range("C28").Value = -18.23663902282715
$ws.Range("D28").Value = -27.17941665649414
$ws.Range("E28").Value = -57.50062561035156
$ws.Range("F28").Value = 2.691667430955088
$ws.Range("G28").Value = 9.822626233541001
$ws.Range("H28").Value = -7.046651205013665

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "walkingToRunning"
$ws.Range("C29").Value = -9.547454833984377
$ws.Range("D29").Value = -6.686841011047363
$ws.Range("E29").Value = -6.594038009643555
$ws.Range("F29").Value = 4.03936266899109
$ws.Range("G29").Value = -6.018842667231255
$ws.Range("H29").Value = -2.086610792307836

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "walkingToRunning"
$ws.Range("C30").Value = -7.443731307983398
$ws.Range("D30").Value = -33.3333854675293
$ws.Range("E30").Value = -12.94254684448242
$ws.Range("F30").Value = 3.831653899372287
$ws.Range("G30").Value = 3.543843676683208
$ws.Range("H30").Value = -1.997865703273085

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "walkingToRunning"
$ws.Range("C31").Value = -11.20403289794922
$ws.Range("D31").Value = -4.358112335205078
$ws.Range("E31").Value = 5.592035293579102
$ws.Range("F31").Value = -3.046866848899856
$ws.Range("G31").Value = 0.7223788134726359
$ws.Range("H31").Value = 2.788422576615722

